$d = $word.ActiveDocument

# 1. Merge "CW" + "2 :" + " Group Collaboration Summary" into a single run's text.
$d.Content.Find.Execute("CW2 : Group Collaboration Summary", $true, $false, $false, $false, $false, $true, 1, $false, "CW2 : Group Collaboration Summary", 2)

# 2. Remove the leading "(" before "The rest of the classes..."
$d.Content.Find.Execute("(The rest of the classes in the application, Junit testing, general application testing.", $true, $false, $false, $false, $false, $true, 1, $false, "The rest of the classes in the application, Junit testing, general application testing.", 2)

# 3. Merge "Performing " and "Overview video." texts (already contiguous, this is a no-op text-wise)
$d.Content.Find.Execute("Performing Overview video.", $true, $false, $false, $false, $false, $true, 1, $false, "Performing Overview video.", 2)

# 4. Merge "Deck class, " + "CardNode" + " class, a share of the Junit testing and general testing of the application."
$d.Content.Find.Execute("Deck class, CardNode class, a share of the Junit testing and general testing of the application.", $true, $false, $false, $false, $false, $true, 1, $false, "Deck class, CardNode class, a share of the Junit testing and general testing of the application.", 2)
